$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 287; this shifts the existing rows 287..410 down to 288..411
# and automatically extends the used range / dimension to row 411.
$ws.Rows(287).Insert()

# Populate the newly inserted row 287 with the new data record.
$ws.Cells.Item(287, 1).Value = 9
$ws.Cells.Item(287, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(287, 3).Value = "Metropolitana"

$ws.Cells.Item(287, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(287, 4).Value = 44726

$ws.Cells.Item(287, 5).Value = 13
$ws.Cells.Item(287, 6).Value = 100112039
$ws.Cells.Item(287, 7).Value = "Ciboulette"
$ws.Cells.Item(287, 8).Value = "Sin especificar"
$ws.Cells.Item(287, 9).Value = "Primera"
$ws.Cells.Item(287, 10).Value = 680
$ws.Cells.Item(287, 11).Value = 800
$ws.Cells.Item(287, 12).Value = 1000
$ws.Cells.Item(287, 13).Value = 912
$ws.Cells.Item(287, 14).Value = '$/docena de atados'
$ws.Cells.Item(287, 15).Value = "Región Metropolitana"
$ws.Cells.Item(287, 16).Value = 304
$ws.Cells.Item(287, 17).Value = 3
$ws.Cells.Item(287, 18).Value = "Hortaliza"
